# Finished updating rules. Started implementing Sherman Movement rules and Firing rules.
#
# Inserts two new event rows ("e052 Pivot Tank" and "e053 Main Gun Firing")
# into the Events sheet, right before the existing "e100 Evening Debriefing"
# row (which was row 62 and is now pushed down to row 64).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

# Insert two new blank rows at 62/63 - everything from the old row 62
# ("e100 Evening Debriefing", previously A62/B62) onward shifts down by 2.
$ws.Rows("62:63").Insert() | Out-Null

# --- Row 62: e052 Pivot Tank ---------------------------------------------
$ws.Range("A62").Value = "e052"
$ws.Range("B62").Value = "<Bold>e052 Pivot Tank</Bold> `n<InlineUIContainer><Button Content='r4.74.1' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   `n<InlineUIContainer><Button Content='r8.46' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `n<LineBreak/><LineBreak/>`nSelect the plus or minus buttons to rotate.  `n<LineBreak/><LineBreak/>`n                                   <InlineUIContainer><Button Content='   -   ' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `n<InlineUIContainer><Image Name='ShermanPivot'  Height='150' Width='150'></Image></InlineUIContainer> `n<InlineUIContainer><Button Content='   +   ' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>`n<LineBreak/><LineBreak/>`nWhen you are satisfied with the current orientation, click Sherman image between buttons to continue."
$ws.Rows.Item(62).RowHeight = 165

# --- Row 63: e053 Main Gun Firing -----------------------------------------
$ws.Range("A63").Value = "e053"
$ws.Range("B63").Value = "<Bold>e053 Main Gun Firing</Bold> `n<InlineUIContainer><Button Content='r4.74.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   `n<InlineUIContainer><Button Content='r9.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `n<LineBreak/><LineBreak/>`nSelect a target by clicking the enemy unit. Consult the <InlineUIContainer><Button Content='To Hit Target' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Table to determine if target is hit.`n<LineBreak/><LineBreak/>`n<Underline>Modifiers:</Underline><LineBreak/>"
$ws.Rows.Item(63).RowHeight = 120

# Match the author's final view state: scrolled down a bit further, with
# B63 (the newly added "e053" description cell) selected.
$ws.Application.ActiveWindow.ScrollRow = 61
$ws.Range("B63").Select() | Out-Null

Write-Output "Inserted e052/e053 rows into Events sheet"
